$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# New report data (post "Generate Report for Handoff" run)
# -----------------------------------------------------------------
$mdFile1Old = "046f5f43-136a-40ec-b0a8-eb093b422b4b.md"
$mdFile2Old = "fe7cc5c2-cea7-482f-8494-90acb7a68223.md"
$mdFile1New = "f64fd7e2-04ab-449d-ada8-fdaffe3197cb.md"
$mdFile2New = "ffffb99fb7df-7bc5-4f5a-a233-b7bb0c6e0563.md"

$status = "Ready for handoff"

$zhHandoffOld = "046f5f43-136a-40ec-b0a8-eb093b422b4b.344f6ed9f063a3ff77d7a5e451e0ce49fc95b772.zh-cn.xlf"
$deHandoffOld = "046f5f43-136a-40ec-b0a8-eb093b422b4b.344f6ed9f063a3ff77d7a5e451e0ce49fc95b772.de-de.xlf"

$zhHandoffNew = "f64fd7e2-04ab-449d-ada8-fdaffe3197cb.c4ef44521985cc8052aa5530f95c3ba80f4971b4.zh-cn.xlf"
$deHandoffNew = "f64fd7e2-04ab-449d-ada8-fdaffe3197cb.c4ef44521985cc8052aa5530f95c3ba80f4971b4.de-de.xlf"

$zhHandoffDatetime = "2016-03-13 21:13:27"
$deHandoffDatetime = "2016-03-13 21:13:33"
$handbackDatetime = "0001-01-01 00:00:00"

# External hyperlink targets (unchanged by this edit - keep the existing ones)
$urlMd1 = "https://github.com/OpenLocalizationTest/oltest/blob/092f90444e72581942c9ed22db27d282a15dd546/e2e/046f5f43-136a-40ec-b0a8-eb093b422b4b.md"
$urlMd2 = "https://github.com/OpenLocalizationTest/oltest/blob/092f90444e72581942c9ed22db27d282a15dd546/e2e/fe7cc5c2-cea7-482f-8494-90acb7a68223.md"

$urlZhHandoff = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1d871c72076c670a54c93d10b43116d12d5ed7af/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/046f5f43-136a-40ec-b0a8-eb093b422b4b.344f6ed9f063a3ff77d7a5e451e0ce49fc95b772.zh-cn.xlf"
$urlDeHandoff = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/94e7763018fab8506b38b71dd549bd9b56694481/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/046f5f43-136a-40ec-b0a8-eb093b422b4b.344f6ed9f063a3ff77d7a5e451e0ce49fc95b772.de-de.xlf"

# -----------------------------------------------------------------
# Sheet "Overview": update the two md-file hyperlinked rows
# -----------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Hyperlinks.Delete()

$wsOverview.Range("A2").Value = $mdFile1New
$wsOverview.Range("A3").Value = $mdFile2New

$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $urlMd1, "", "", $mdFile1New)
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), $urlMd2, "", "", $mdFile2New)

# -----------------------------------------------------------------
# Sheet "zh-cn"
# -----------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Hyperlinks.Delete()

# Row 2 (file 1)
$wsZh.Range("A2").Value = $mdFile1New
$wsZh.Range("B2").Value = ".md"
$wsZh.Range("C2").Value = $status
$wsZh.Range("D2").Value = $zhHandoffNew
$wsZh.Range("E2").Value = $zhHandoffDatetime
$wsZh.Range("F2").Clear()
$wsZh.Range("G2").Clear()
$wsZh.Range("H2").Value = $handbackDatetime
$wsZh.Range("I2").Value = "Include"

# Row 3 (file 2)
$wsZh.Range("A3").Value = $mdFile2New
$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = $status
$wsZh.Range("D3").Value = $zhHandoffNew
$wsZh.Range("E3").Value = $zhHandoffDatetime
$wsZh.Range("F3").Clear()
$wsZh.Range("G3").Clear()
$wsZh.Range("H3").Value = $handbackDatetime
$wsZh.Range("I3").Value = "Include"

$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $urlMd1, "", "", $mdFile1New)
$wsZh.Hyperlinks.Add($wsZh.Range("B2"), $urlMd1, "", "", ".md")
$wsZh.Hyperlinks.Add($wsZh.Range("D2"), $urlZhHandoff, "", "", $zhHandoffNew)
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $urlMd2, "", "", $mdFile2New)
$wsZh.Hyperlinks.Add($wsZh.Range("B3"), $urlMd2, "", "", ".md")
$wsZh.Hyperlinks.Add($wsZh.Range("D3"), $urlZhHandoff, "", "", $zhHandoffNew)

# -----------------------------------------------------------------
# Sheet "de-de"
# -----------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Hyperlinks.Delete()

# Row 2 (file 1)
$wsDe.Range("A2").Value = $mdFile1New
$wsDe.Range("B2").Value = ".md"
$wsDe.Range("C2").Value = $status
$wsDe.Range("D2").Value = $deHandoffNew
$wsDe.Range("E2").Value = $deHandoffDatetime
$wsDe.Range("F2").Clear()
$wsDe.Range("G2").Clear()
$wsDe.Range("H2").Value = $handbackDatetime
$wsDe.Range("I2").Value = "Include"

# Row 3 (file 2)
$wsDe.Range("A3").Value = $mdFile2New
$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = $status
$wsDe.Range("D3").Value = $deHandoffNew
$wsDe.Range("E3").Value = $deHandoffDatetime
$wsDe.Range("F3").Clear()
$wsDe.Range("G3").Clear()
$wsDe.Range("H3").Value = $handbackDatetime
$wsDe.Range("I3").Value = "Include"

$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $urlMd1, "", "", $mdFile1New)
$wsDe.Hyperlinks.Add($wsDe.Range("B2"), $urlMd1, "", "", ".md")
$wsDe.Hyperlinks.Add($wsDe.Range("D2"), $urlDeHandoff, "", "", $deHandoffNew)
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $urlMd2, "", "", $mdFile2New)
$wsDe.Hyperlinks.Add($wsDe.Range("B3"), $urlMd2, "", "", ".md")
$wsDe.Hyperlinks.Add($wsDe.Range("D3"), $urlDeHandoff, "", "", $deHandoffNew)
